# Weekly update: a new "Terminal Hortofrutícola Agro Chillán - Cilantro" price
# record (for market date serial 44875) is inserted at the top of the data
# block (row 92), pushing all existing price rows down by two rows.
#
# Columns: A Mercado ID | B Mercado | C Región | D Fecha | E Codreg |
#          F Categoría ID | G Categoría | H Variedad | I Calidad | J Volumen |
#          K Precio mínimo | L Precio máximo | M Precio promedio ponderado |
#          N Unidad de comercialización | O Origen | P Precio $/Kg |
#          Q Kg o Unidades | R Clasificación

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows above the current row 92, shifting the existing
# rows 92-113 down to rows 94-115 (dimension grows from A1:R113 to A1:R115).
$ws.Range("A92:A93").EntireRow.Insert()

# Shared / constant values for this product across every row in the block.
$mercadoId  = 7
$mercado    = "Terminal Hortofrutícola Agro Chillán"
$region     = "Ñuble"
$codreg     = 16
$categoriaId = 100112040
$categoria  = "Cilantro"
$variedad   = "Sin especificar"
$unidad     = "$/atado 0,5 a 1 kilo"
$origen     = "Provincia de Diguillín"
$kgOUnidades = 1
$clasificacion = "Hortaliza"
$fecha = 44875

# --- New row 92: Calidad "Primera" ---
$r = 92
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $categoriaId
$ws.Cells.Item($r, 7).Value  = $categoria
$ws.Cells.Item($r, 8).Value  = $variedad
$ws.Cells.Item($r, 9).Value  = "Primera"
$ws.Cells.Item($r, 10).Value = 300
$ws.Cells.Item($r, 11).Value = 600
$ws.Cells.Item($r, 12).Value = 700
$ws.Cells.Item($r, 13).Value = 650
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = $origen
$ws.Cells.Item($r, 16).Value = 650
$ws.Cells.Item($r, 17).Value = $kgOUnidades
$ws.Cells.Item($r, 18).Value = $clasificacion

# --- New row 93: Calidad "Segunda" ---
$r = 93
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $categoriaId
$ws.Cells.Item($r, 7).Value  = $categoria
$ws.Cells.Item($r, 8).Value  = $variedad
$ws.Cells.Item($r, 9).Value  = "Segunda"
$ws.Cells.Item($r, 10).Value = 200
$ws.Cells.Item($r, 11).Value = 500
$ws.Cells.Item($r, 12).Value = 500
$ws.Cells.Item($r, 13).Value = 500
$ws.Cells.Item($r, 14).Value = $unidad
$ws.Cells.Item($r, 15).Value = $origen
$ws.Cells.Item($r, 16).Value = 500
$ws.Cells.Item($r, 17).Value = $kgOUnidades
$ws.Cells.Item($r, 18).Value = $clasificacion
